$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty row 43, shifting rows 44:67 up by one.
$ws.Rows.Item(43).Delete()

# Update the selected cell / view as in the target workbook.
$ws.Range("H13").Select()
